{"js": "// AUTO FROM WORK 10.05.2022 11:49:20,45\n// 1) Split the \"21 24\" run into separate \"21\" and \"24\" runs, highlight\n//    each in yellow, and relocate the \"_GoBack\" bookmark so it still\n//    marks the end of the edited text (Word keeps it pinned to the most\n//    recently edited spot).\nconst body = context.document.body;\n\nconst combinedResults = body.search(\"21 24\", { matchCase: true });\ncombinedResults.load(\"items\");\nawait context.sync();\n\nconst combined = combinedResults.items[0];\n\nconst run21Results = combined.search(\"21\", { matchCase: true });\nconst run24Results = combined.search(\"24\", { matchCase: true });\nrun21Results.load(\"items\");\nrun24Results.load(\"items\");\nawait context.sync();\n\nconst run21 = run21Results.items[0];\nconst run24 = run24Results.items[0];\nrun21.font.highlightColor = \"Yellow\";\nrun24.font.highlightColor = \"Yellow\";\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst afterRun24 = run24.getRange(\"End\");\nafterRun24.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Highlight the \"\u041d\u0430\u0439\u0442\u0438 \u0441\u0443\u043c\u043c\u0443 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u043e\u0432 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u043c\u0435\u0436\u0434\u0443\n//    \u043f\u0435\u0440\u0432\u044b\u043c \u0438 \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043d\u0443\u043b\u0435\u0432\u044b\u043c\u0438 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u0430\u043c\u0438.\" bullet (paragraph mark +\n//    both runs).\nconst para2Results = body.search(\n  \"\u041d\u0430\u0439\u0442\u0438 \u0441\u0443\u043c\u043c\u0443 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u043e\u0432 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u043c\u0435\u0436\u0434\u0443 \u043f\u0435\u0440\u0432\u044b\u043c \u0438 \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043d\u0443\u043b\u0435\u0432\u044b\u043c\u0438 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u0430\u043c\u0438.\",\n  { matchCase: true }\n);\npara2Results.load(\"items\");\nawait context.sync();\n\nconst para2Match = para2Results.items[0];\nconst para2 = para2Match.paragraphs.getFirst();\nconst para2Range = para2.getRange();\npara2.font.highlightColor = \"Yellow\";\npara2Range.font.highlightColor = \"Yellow\";\nawait context.sync();\n\n// 3) Highlight the \"\u041d\u0430\u043f\u0438\u0441\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0443, \u0432 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0438 \u0437\u0430\u043f\u043e\u043b\u043d\u0438\u0442\u044c\n//    \u0441\u043b\u0443\u0447\u0430\u0439\u043d\u044b\u043c\u0438 \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u044f\u043c\u0438 \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435 0..9 ... \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u0432\u044b\u0448\u0435\n//    \u0433\u043b\u0430\u0432\u043d\u043e\u0439 \u0434\u0438\u0430\u0433\u043e\u043d\u0430\u043b\u0438\" bullet (paragraph mark + all runs).\nconst para3Results = body.search(\n  \"\u041d\u0430\u043f\u0438\u0441\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0443, \u0432 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0438 \u0437\u0430\u043f\u043e\u043b\u043d\u0438\u0442\u044c \u0441\u043b\u0443\u0447\u0430\u0439\u043d\u044b\u043c\u0438 \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u044f\u043c\u0438 \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435\",\n  { matchCase: true }\n);\npara3Results.load(\"items\");\nawait context.sync();\n\nconst para3Match = para3Results.items[0];\nconst para3 = para3Match.paragraphs.getFirst();\nconst para3Range = para3.getRange();\npara3.font.highlightColor = \"Yellow\";\npara3Range.font.highlightColor = \"Yellow\";\nawait context.sync();\n", "ps1": "# AUTO FROM WORK 10.05.2022 11:49:20,45\n$d = $word.ActiveDocument\n\n# 1) Split the \"21 24\" run into separate \"21\" and \"24\" pieces, highlight\n#    each in yellow (the space between them stays un-highlighted), and\n#    relocate the \"_GoBack\" bookmark so it still marks the end of the\n#    edited text (Word keeps it pinned to the most recently edited spot).\n$combined = $d.Content\n$found = $combined.Find.Execute(\"21 24\")\n$startPos = $combined.Start\n$endPos = $combined.End\n\n$run21 = $d.Range($startPos, $startPos)\n$run21.MoveEnd(\"wdCharacter\", 2)\n$run21.Font.HighlightColorIndex = \"wdYellow\"\n\n$run24 = $d.Range($endPos, $endPos)\n$run24.MoveStart(\"wdCharacter\", -2)\n$run24.Font.HighlightColorIndex = \"wdYellow\"\n\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n$goBackTarget = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $goBackTarget)\n\n# 2) Highlight the \"\u041d\u0430\u0439\u0442\u0438 \u0441\u0443\u043c\u043c\u0443 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u043e\u0432 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u043c\u0435\u0436\u0434\u0443\n#    \u043f\u0435\u0440\u0432\u044b\u043c \u0438 \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043d\u0443\u043b\u0435\u0432\u044b\u043c\u0438 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u0430\u043c\u0438.\" bullet (paragraph mark +\n#    both runs).\n$find2 = $d.Content\n$find2.Find.Execute(\"\u041d\u0430\u0439\u0442\u0438 \u0441\u0443\u043c\u043c\u0443 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u043e\u0432 \u043c\u0430\u0441\u0441\u0438\u0432\u0430, \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u043c\u0435\u0436\u0434\u0443 \u043f\u0435\u0440\u0432\u044b\u043c \u0438 \u043f\u043e\u0441\u043b\u0435\u0434\u043d\u0438\u043c \u043d\u0443\u043b\u0435\u0432\u044b\u043c\u0438 \u044d\u043b\u0435\u043c\u0435\u043d\u0442\u0430\u043c\u0438.\")\n$para2Range = $find2.Paragraphs(1).Range\n$para2Range.Font.HighlightColorIndex = \"wdYellow\"\n\n# 3) Highlight the \"\u041d\u0430\u043f\u0438\u0441\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0443, \u0432 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0438 \u0437\u0430\u043f\u043e\u043b\u043d\u0438\u0442\u044c\n#    \u0441\u043b\u0443\u0447\u0430\u0439\u043d\u044b\u043c\u0438 \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u044f\u043c\u0438 \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435 0..9 ... \u0440\u0430\u0441\u043f\u043e\u043b\u043e\u0436\u0435\u043d\u043d\u044b\u0445 \u0432\u044b\u0448\u0435\n#    \u0433\u043b\u0430\u0432\u043d\u043e\u0439 \u0434\u0438\u0430\u0433\u043e\u043d\u0430\u043b\u0438\" bullet (paragraph mark + all runs).\n$find3 = $d.Content\n$find3.Find.Execute(\"\u041d\u0430\u043f\u0438\u0441\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u0440\u0430\u043c\u043c\u0443, \u0432 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u0441\u043e\u0437\u0434\u0430\u0442\u044c \u0438 \u0437\u0430\u043f\u043e\u043b\u043d\u0438\u0442\u044c \u0441\u043b\u0443\u0447\u0430\u0439\u043d\u044b\u043c\u0438 \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u044f\u043c\u0438 \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435\")\n$para3Range = $find3.Paragraphs(1).Range\n$para3Range.Font.HighlightColorIndex = \"wdYellow\"\n"}
